# "list out all sheets"
#
# The "Name" column on "Sample xsheet C" gets updated from the placeholder
# John/Sam/Dean rows to Chris/"Chris Wong"/chriswong, the column is
# re-sized to fit the new text, and that sheet becomes the active
# (selected) sheet of the workbook instead of "Sample xsheet A".

$wb = $excel.ActiveWorkbook

$wsC = $wb.Worksheets.Item("Sample xsheet C")

$wsC.Range("B2").Value = "Chris"
$wsC.Range("B3").Value = "Chris Wong"
$wsC.Range("B4").Value = "chriswong"

# Fit column B to the new, longer values (e.g. "chriswong").
[void]$wsC.Columns.Item(2).AutoFit()

# Make "Sample xsheet C" the active sheet / tab, with the whole sheet
# selected - this replaces "Sample xsheet A" as the active tab.
[void]$wsC.Select()
[void]$wsC.Cells.Select()
